$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.862146333333333
$ws.Range("H2").Value = 8.586439
$ws.Range("I2").Value = 0.9778268096017091
$ws.Range("J2").Value = 0.9778268096017091
$ws.Range("M2").Value = 0.7564453333333333
$ws.Range("N2").Value = 2.269336
$ws.Range("O2").Value = 0.2228868660598179
$ws.Range("P2").Value = 0.2228868660598179
$ws.Range("Q2").Value = 2.165057237167111
$ws.Range("R2").Value = 19.485515134504
$ws.Range("S2").Value = 0.2179447531413952
$ws.Range("T2").Value = 0.2179447531413952
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.862146333333333
$ws.Range("H3").Value = 8.586439
$ws.Range("I3").Value = 0.9778268096017091
$ws.Range("J3").Value = 0.9778268096017091
$ws.Range("O3").Value = 0.5246945458259835
$ws.Range("P3").Value = 0.5246945458259834
$ws.Range("Q3").Value = 5.096727967083445
$ws.Range("R3").Value = 45.870551703751
$ws.Range("S3").Value = 0.5130603937604392
$ws.Range("T3").Value = 0.5130603937604391
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.862146333333333
$ws.Range("H4").Value = 8.586439
$ws.Range("I4").Value = 0.9778268096017091
$ws.Range("J4").Value = 0.9778268096017091
$ws.Range("M4").Value = 0.6826836666666667
$ws.Range("N4").Value = 2.048051
$ws.Range("O4").Value = 0.2011529667359422
$ws.Range("P4").Value = 0.2011529667359422
$ws.Range("Q4").Value = 1.953940553376556
$ws.Range("R4").Value = 17.585464980389
$ws.Range("S4").Value = 0.1966927637053251
$ws.Range("T4").Value = 0.1966927637053251
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.862146333333333
$ws.Range("H5").Value = 8.586439
$ws.Range("I5").Value = 0.9778268096017091
$ws.Range("J5").Value = 0.9778268096017091
$ws.Range("M5").Value = 0.173988
$ws.Range("N5").Value = 0.5219640000000001
$ws.Range("O5").Value = 0.05126562137825639
$ws.Range("P5").Value = 0.05126562137825638
$ws.Range("Q5").Value = 0.4979791162440001
$ws.Range("R5").Value = 4.481812046196001
$ws.Range("S5").Value = 0.05012889899454961
$ws.Range("T5").Value = 0.05012889899454961
$ws.Range("I6").Value = 0.02217319039829088
$ws.Range("J6").Value = 0.02217319039829088
$ws.Range("M6").Value = 0.7564453333333333
$ws.Range("N6").Value = 2.269336
$ws.Range("O6").Value = 0.2228868660598179
$ws.Range("P6").Value = 0.2228868660598179
$ws.Range("Q6").Value = 0.049094815024
$ws.Range("R6").Value = 0.441853335216
$ws.Range("S6").Value = 0.004942112918422701
$ws.Range("T6").Value = 0.0049421129184227
$ws.Range("I7").Value = 0.02217319039829088
$ws.Range("J7").Value = 0.02217319039829088
$ws.Range("O7").Value = 0.5246945458259835
$ws.Range("P7").Value = 0.5246945458259834
$ws.Range("S7").Value = 0.01163415206554429
$ws.Range("T7").Value = 0.01163415206554429
$ws.Range("I8").Value = 0.02217319039829088
$ws.Range("J8").Value = 0.02217319039829088
$ws.Range("M8").Value = 0.6826836666666667
$ws.Range("N8").Value = 2.048051
$ws.Range("O8").Value = 0.2011529667359422
$ws.Range("P8").Value = 0.2011529667359422
$ws.Range("Q8").Value = 0.044307535334
$ws.Range("R8").Value = 0.398767818006
$ws.Range("S8").Value = 0.004460203030617119
$ws.Range("T8").Value = 0.004460203030617119
$ws.Range("I9").Value = 0.02217319039829088
$ws.Range("J9").Value = 0.02217319039829088
$ws.Range("M9").Value = 0.173988
$ws.Range("N9").Value = 0.5219640000000001
$ws.Range("O9").Value = 0.05126562137825639
$ws.Range("P9").Value = 0.05126562137825638
$ws.Range("Q9").Value = 0.011292169176
$ws.Range("R9").Value = 0.101629522584
$ws.Range("S9").Value = 0.00113672238370677
$ws.Range("T9").Value = 0.00113672238370677

Write-Host "done"